$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# Remove the redundant run "    " (4 spaces) that sits right between
# "...NotExistingRepo" and "demonstration" in the first paragraph. There is
# a similar looking run earlier in the paragraph (between "A simple" and
# "<---") which must stay untouched, so locate this occurrence precisely
# via Find on the surrounding text before deleting just the space span.
$rng = $d.Content
$found = $rng.Find.Execute("NotExistingRepo    demonstration", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $spacesStart = $rng.Start + 15
    $spacesEnd = $spacesStart + 4
    $spacesRng = $d.Range($spacesStart, $spacesEnd)
    $spacesRng.Delete()
}

# --- Change 2 -----------------------------------------------------------
# Replace the " m:self.name " field (fldChar begin / instrText* / fldChar
# end) with plain-text runs spelling out "{m:self.name}" -- i.e. drop the
# field wrapper and the stray spaces, keep the "self" run's color
# formatting, and close the brace directly after "name".
$f = $d.Fields(1)
$f.Delete()

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>self</w:t></w:r><w:r><w:t xml:space="preserve">.name}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$fieldParaRange = $d.Paragraphs(2).Range
[void]$fieldParaRange.InsertXML($xml)
